# Scheduled market-price refresh: updates currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns (H:N) on the Leve-profit
# sheets for the rows where fresh marketboard data came back. Columns
# A:G (leve/item identity) are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1531.5385
$ws.Range("I15").Value = 1531.5385
$ws.Range("K15").Value = 4594.6155
$ws.Range("M15").Value = -4425.6155

$ws.Range("H51").Value = 8337281
$ws.Range("I51").Value = 3995
$ws.Range("K51").Value = 3995
$ws.Range("M51").Value = -3511

$ws.Range("H53").Value = 71.818184
$ws.Range("I53").Value = 45.833332
$ws.Range("J53").Value = 103
$ws.Range("K53").Value = 45.833332
$ws.Range("L53").Value = 103
$ws.Range("M53").Value = 591.166668
$ws.Range("N53").Value = -1377

$ws.Range("H86").Value = 1300

$ws.Range("H89").Value = 1300

$ws.Range("H132").Value = 1419.6809
$ws.Range("I132").Value = 1419.6809
$ws.Range("K132").Value = 4259.0427
$ws.Range("M132").Value = -1729.0427

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5749.974
$ws.Range("I32").Value = 2442.4307
$ws.Range("J32").Value = 53378.6
$ws.Range("K32").Value = 2442.4307
$ws.Range("L32").Value = 53378.6
$ws.Range("M32").Value = -2155.4307
$ws.Range("N32").Value = -53952.6

$ws.Range("H61").Value = 1996.2
$ws.Range("I61").Value = 2128.6667
$ws.Range("J61").Value = 1797.5
$ws.Range("K61").Value = 2128.6667
$ws.Range("L61").Value = 1797.5
$ws.Range("M61").Value = -1916.6667
$ws.Range("N61").Value = -2221.5

$ws.Range("H97").Value = 1083.9354
$ws.Range("I97").Value = 1156.0741
$ws.Range("K97").Value = 1156.0741
$ws.Range("M97").Value = -660.0741

$ws.Range("H136").Value = 1996.2
$ws.Range("I136").Value = 2128.6667
$ws.Range("J136").Value = 1797.5
$ws.Range("K136").Value = 6386.000100000001
$ws.Range("L136").Value = 5392.5
$ws.Range("M136").Value = -3836.000100000001
$ws.Range("N136").Value = -10492.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 22777.428
$ws.Range("I76").Value = 23000
$ws.Range("K76").Value = 23000
$ws.Range("M76").Value = -22685

$ws.Range("H79").Value = 22777.428
$ws.Range("I79").Value = 23000
$ws.Range("K79").Value = 23000
$ws.Range("M79").Value = -21908

$ws.Range("H86").Value = 1510.4706
$ws.Range("I86").Value = 1237.5
$ws.Range("K86").Value = 1237.5
$ws.Range("M86").Value = -114.5

$ws.Range("H89").Value = 1510.4706
$ws.Range("I89").Value = 1237.5
$ws.Range("K89").Value = 6187.5
$ws.Range("M89").Value = -571.5

$ws.Range("H94").Value = 1960.5
$ws.Range("I94").Value = 2052.6
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 2052.6
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -1601.6
$ws.Range("N94").Value = -2402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31543.521
$ws.Range("I31").Value = 29520.568
$ws.Range("J31").Value = 39860.11
$ws.Range("K31").Value = 29520.568
$ws.Range("L31").Value = 39860.11
$ws.Range("M31").Value = -29225.568
$ws.Range("N31").Value = -40450.11

$ws.Range("H34").Value = 31543.521
$ws.Range("I34").Value = 29520.568
$ws.Range("J34").Value = 39860.11
$ws.Range("K34").Value = 29520.568
$ws.Range("L34").Value = 39860.11
$ws.Range("M34").Value = -29318.568
$ws.Range("N34").Value = -40264.11

$ws.Range("H86").Value = 3636.2856
$ws.Range("I86").Value = 3366.25
$ws.Range("K86").Value = 3366.25
$ws.Range("M86").Value = -2243.25

$ws.Range("H89").Value = 3636.2856
$ws.Range("I89").Value = 3366.25
$ws.Range("K89").Value = 16831.25
$ws.Range("M89").Value = -11215.25

$ws.Range("H105").Value = 1457.5
$ws.Range("I105").Value = 909.5
$ws.Range("K105").Value = 909.5
$ws.Range("M105").Value = 837.5

$ws.Range("H134").Value = 7649.5713
$ws.Range("I134").Value = 5886.75
$ws.Range("K134").Value = 17660.25
$ws.Range("M134").Value = -15125.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 562.2632
$ws.Range("J86").Value = 717
$ws.Range("L86").Value = 2151
$ws.Range("N86").Value = -4523

$ws.Range("H89").Value = 562.2632
$ws.Range("J89").Value = 717
$ws.Range("L89").Value = 6453
$ws.Range("N89").Value = -18309

$ws.Range("H113").Value = 2949.3333
$ws.Range("J113").Value = 3999.5
$ws.Range("L113").Value = 11998.5
$ws.Range("N113").Value = -16338.5

$ws.Range("H122").Value = 875.1429000000001
$ws.Range("J122").Value = 905.4
$ws.Range("L122").Value = 8148.599999999999
$ws.Range("N122").Value = -13048.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3800.6
$ws.Range("I80").Value = 3126.7273
$ws.Range("J80").Value = 4624.222
$ws.Range("K80").Value = 3126.7273
$ws.Range("L80").Value = 4624.222
$ws.Range("M80").Value = -2128.7273
$ws.Range("N80").Value = -6620.222

$ws.Range("H83").Value = 3800.6
$ws.Range("I83").Value = 3126.7273
$ws.Range("J83").Value = 4624.222
$ws.Range("K83").Value = 15633.6365
$ws.Range("L83").Value = 23121.11
$ws.Range("M83").Value = -10641.6365
$ws.Range("N83").Value = -33105.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H68").Value = 3430
$ws.Range("I68").Value = 3200
$ws.Range("J68").Value = 3545
$ws.Range("K68").Value = 3200
$ws.Range("L68").Value = 3545
$ws.Range("M68").Value = -2451
$ws.Range("N68").Value = -5043

$ws.Range("H71").Value = 3430
$ws.Range("I71").Value = 3200
$ws.Range("J71").Value = 3545
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 17725
$ws.Range("M71").Value = -12256
$ws.Range("N71").Value = -25213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20755
$ws.Range("I62").Value = 5050.3335
$ws.Range("J62").Value = 39600.6
$ws.Range("K62").Value = 5050.3335
$ws.Range("L62").Value = 39600.6
$ws.Range("M62").Value = -4426.3335
$ws.Range("N62").Value = -40848.6

$ws.Range("H65").Value = 20755
$ws.Range("I65").Value = 5050.3335
$ws.Range("J65").Value = 39600.6
$ws.Range("K65").Value = 25251.6675
$ws.Range("L65").Value = 198003
$ws.Range("M65").Value = -22131.6675
$ws.Range("N65").Value = -204243

$ws.Range("H96").Value = 2963.1924
$ws.Range("I96").Value = 1398.25
$ws.Range("K96").Value = 1398.25
$ws.Range("M96").Value = -25.25

$ws.Range("H100").Value = 2476.2307
$ws.Range("I100").Value = 2455.3333
$ws.Range("K100").Value = 4910.6666
$ws.Range("M100").Value = -4369.6666
